$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - 100m Male results
$ws.Range("B2").Value = "Noah Lyles"
$ws.Range("C2").Value = "Kishane Thompson"
$ws.Range("D2").Value = "Fred Kerley"

# Row 3 - 100m Female results
$ws.Range("B3").Value = "Julien Alfred"
$ws.Range("D3").Value = "Melissa Jefferson"
$ws.Range("C3").Value = "Sha'Carri Richardson"

# Autofit columns B:D to match bestFit/auto-sizing behavior after data entry
$ws.Range("B1:D9").EntireColumn.AutoFit() | Out-Null

# Nudge the autofit widths to line up with Excel's own font-metric based
# best-fit calculation for these specific values
$ws.Range("B1").ColumnWidth = 11.7109375
$ws.Range("C1").ColumnWidth = 19.85546875
$ws.Range("D1").ColumnWidth = 16.28515625

# Update selection to match the recorded cursor position
$ws.Range("H5").Select() | Out-Null
